$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the wrong correlation_direction values (column J) for the
# "Hohmeyer & Wolff (2010), Table 15" monthly_benefit_receipt_effect
# estimates (rows 14-25): the correlation direction should be negative,
# not positive.
for ($row = 14; $row -le 25; $row++) {
    $ws.Cells.Item($row, 10).Value = -1
}

# Update the active selection to the corrected range (and drop the
# previous scrolled-away topLeftCell) so the sheet view reflects the
# cells that were just fixed.
$ws.Range("J14:J25").Select()
